$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts catatan/nilai/tgl_mitra_diterima right)
$ws.Columns("D:D").Insert()

# Populate the new "rate_honor" column header and example value
$ws.Range("D1").Value = "rate_honor"
$ws.Range("D2").Value = "1"

# Update the selected cell shown in the sheet view
[void]$ws.Range("G10").Select()
